# Auto-generated script to update Halicarnassus_Profits market-price cells
# across all item-category worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 518.25
$ws.Range("I6").Value = 340
$ws.Range("K6").Value = 1020
$ws.Range("M6").Value = -908
$ws.Range("H55").Value = 650.375
$ws.Range("I55").Value = 322.8
$ws.Range("J55").Value = 1196.3334
$ws.Range("K55").Value = 322.8
$ws.Range("L55").Value = 1196.3334
$ws.Range("M55").Value = -108.8
$ws.Range("N55").Value = -1624.3334
$ws.Range("H76").Value = 3019.4
$ws.Range("I76").Value = 3999
$ws.Range("J76").Value = 2366.3333
$ws.Range("K76").Value = 3999
$ws.Range("L76").Value = 2366.3333
$ws.Range("M76").Value = -3684
$ws.Range("N76").Value = -2996.3333
$ws.Range("H79").Value = 3019.4
$ws.Range("I79").Value = 3999
$ws.Range("J79").Value = 2366.3333
$ws.Range("K79").Value = 3999
$ws.Range("L79").Value = 2366.3333
$ws.Range("M79").Value = -2907
$ws.Range("N79").Value = -4550.3333
$ws.Range("H137").Value = 15000
$ws.Range("J137").Value = 15000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -50100
$ws.Range("H138").Value = 4620
$ws.Range("J138").Value = 4892
$ws.Range("L138").Value = 14676
$ws.Range("N138").Value = -24956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 896.4
$ws.Range("I2").Value = 840.44446
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 840.44446
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -727.44446
$ws.Range("N2").Value = -1626
$ws.Range("H6").Value = 22500000
$ws.Range("I6").Value = 24000000
$ws.Range("K6").Value = 24000000
$ws.Range("M6").Value = -23999827
$ws.Range("H32").Value = 3318.6775
$ws.Range("I32").Value = 3202.7585
$ws.Range("K32").Value = 3202.7585
$ws.Range("M32").Value = -2915.7585
$ws.Range("H36").Value = 2600
$ws.Range("I36").Value = 2600
$ws.Range("K36").Value = 2600
$ws.Range("M36").Value = -2254
$ws.Range("H74").Value = 2338.55
$ws.Range("I74").Value = 1935.3158
$ws.Range("K74").Value = 1935.3158
$ws.Range("M74").Value = -1061.3158
$ws.Range("H77").Value = 2338.55
$ws.Range("I77").Value = 1935.3158
$ws.Range("K77").Value = 9676.579
$ws.Range("M77").Value = -5308.579
$ws.Range("H116").Value = 896.4
$ws.Range("I116").Value = 840.44446
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 840.44446
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 1453.55554
$ws.Range("N116").Value = -5988
$ws.Range("H135").Value = 56085.4
$ws.Range("J135").Value = 56085.4
$ws.Range("L135").Value = 56085.4
$ws.Range("N135").Value = -66225.39999999999
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 896.4
$ws.Range("I3").Value = 840.44446
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 840.44446
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -726.44446
$ws.Range("N3").Value = -1628
$ws.Range("H86").Value = 4903.8184
$ws.Range("I86").Value = 4394.2
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 4394.2
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -3271.2
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 4903.8184
$ws.Range("I89").Value = 4394.2
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 21971
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -16355
$ws.Range("N89").Value = -61232
$ws.Range("H94").Value = 1609.5
$ws.Range("I94").Value = 1609.5
$ws.Range("K94").Value = 1609.5
$ws.Range("M94").Value = -1158.5
$ws.Range("H105").Value = 1916.6666
$ws.Range("I105").Value = 1625
$ws.Range("K105").Value = 1625
$ws.Range("M105").Value = 122

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 315
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("H141").Value = 68563
$ws.Range("J141").Value = 68563
$ws.Range("L141").Value = 68563
$ws.Range("N141").Value = -78923
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1396.2858
$ws.Range("J5").Value = 1621.1428
$ws.Range("L5").Value = 4863.428400000001
$ws.Range("N5").Value = -5087.428400000001
$ws.Range("H17").Value = 2979.8333
$ws.Range("J17").Value = 3816.9285
$ws.Range("L17").Value = 11450.7855
$ws.Range("N17").Value = -11788.7855
$ws.Range("H23").Value = 224.25
$ws.Range("I23").Value = 199
$ws.Range("K23").Value = 597
$ws.Range("M23").Value = -362
$ws.Range("H135").Value = 1396.2858
$ws.Range("J135").Value = 1621.1428
$ws.Range("L135").Value = 14590.2852
$ws.Range("N135").Value = -19660.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 14842785
$ws.Range("I3").Value = 17071090
$ws.Range("J3").Value = 6672333.5
$ws.Range("K3").Value = 17071090
$ws.Range("L3").Value = 6672333.5
$ws.Range("M3").Value = -17070974
$ws.Range("N3").Value = -6672565.5
$ws.Range("H29").Value = 49996.668
$ws.Range("J29").Value = 49996.668
$ws.Range("L29").Value = 49996.668
$ws.Range("N29").Value = -50576.668
$ws.Range("H122").Value = 2733
$ws.Range("I122").Value = 2110.8
$ws.Range("K122").Value = 6332.400000000001
$ws.Range("M122").Value = -3882.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 63552.062
$ws.Range("I20").Value = 1131.0714
$ws.Range("K20").Value = 1131.0714
$ws.Range("M20").Value = -905.0714
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H46").Value = 4430.8
$ws.Range("I46").Value = 2386.3333
$ws.Range("J46").Value = 7497.5
$ws.Range("K46").Value = 2386.3333
$ws.Range("L46").Value = 2386.3333
$ws.Range("M46").Value = -2198.3333
$ws.Range("N46").Value = -7873.5
$ws.Range("H100").Value = 8284.286
$ws.Range("I100").Value = 3995
$ws.Range("K100").Value = 3995
$ws.Range("M100").Value = -3454
$ws.Range("H132").Value = 3133.3333
$ws.Range("I132").Value = 3200
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9600
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7070
$ws.Range("N132").Value = -14060
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("H51").Value = 18520
$ws.Range("I51").Value = 18520
$ws.Range("K51").Value = 18520
$ws.Range("M51").Value = -18010
$ws.Range("N22").Value = -5586
